$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").ClearContents()
$wb.Application.Iterative = $false
